# Roll the SGL 3D-Steel production forecast window forward.
# The sheet holds a rolling 169-row forecast (Data/Interval/Prediction/Lookup)
# keyed by date (A) + hourly interval (B). This refreshes rows 2:170 with the
# latest forecast run: the oldest days drop off the front and new days/
# corrected predictions are appended, while row/column layout, styles and
# headers (row 1) stay untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A = Data (date serial), B = Interval (hour), C = Prediction, D = Lookup (text)
$data = @(
    @(45993, 13, 0.034, "02.12.202513"),
    @(45993, 14, 0.034, "02.12.202514"),
    @(45993, 15, 0.408, "02.12.202515"),
    @(45993, 16, 0.47, "02.12.202516"),
    @(45993, 17, 0.251, "02.12.202517"),
    @(45993, 18, 0.09, "02.12.202518"),
    @(45993, 19, 0.038, "02.12.202519"),
    @(45993, 20, 0.043, "02.12.202520"),
    @(45993, 21, 0.043, "02.12.202521"),
    @(45993, 22, 0.011, "02.12.202522"),
    @(45993, 23, 0.011, "02.12.202523"),
    @(45993, 24, 0.011, "02.12.202524"),
    @(45994, 1, 0, "03.12.20251"),
    @(45994, 2, 0, "03.12.20252"),
    @(45994, 3, 0, "03.12.20253"),
    @(45994, 4, 0, "03.12.20254"),
    @(45994, 5, 0, "03.12.20255"),
    @(45994, 6, 0, "03.12.20256"),
    @(45994, 7, 0, "03.12.20257"),
    @(45994, 8, 0.123, "03.12.20258"),
    @(45994, 9, 0.173, "03.12.20259"),
    @(45994, 10, 0.265, "03.12.202510"),
    @(45994, 11, 0.483, "03.12.202511"),
    @(45994, 12, 0.647, "03.12.202512"),
    @(45994, 13, 0.722, "03.12.202513"),
    @(45994, 14, 0.984, "03.12.202514"),
    @(45994, 15, 0.706, "03.12.202515"),
    @(45994, 16, 0.47, "03.12.202516"),
    @(45994, 17, 0.246, "03.12.202517"),
    @(45994, 18, 0.08599999999999999, "03.12.202518"),
    @(45994, 19, 0.038, "03.12.202519"),
    @(45994, 20, 0.043, "03.12.202520"),
    @(45994, 21, 0.043, "03.12.202521"),
    @(45994, 22, 0.011, "03.12.202522"),
    @(45994, 23, 0.011, "03.12.202523"),
    @(45994, 24, 0.011, "03.12.202524"),
    @(45995, 1, 0, "04.12.20251"),
    @(45995, 2, 0, "04.12.20252"),
    @(45995, 3, 0, "04.12.20253"),
    @(45995, 4, 0, "04.12.20254"),
    @(45995, 5, 0, "04.12.20255"),
    @(45995, 6, 0, "04.12.20256"),
    @(45995, 7, 0, "04.12.20257"),
    @(45995, 8, 0.192, "04.12.20258"),
    @(45995, 9, 0.242, "04.12.20259"),
    @(45995, 10, 0.484, "04.12.202510"),
    @(45995, 11, 0.609, "04.12.202511"),
    @(45995, 12, 0.465, "04.12.202512"),
    @(45995, 13, 0.6879999999999999, "04.12.202513"),
    @(45995, 14, 1.023, "04.12.202514"),
    @(45995, 15, 0.897, "04.12.202515"),
    @(45995, 16, 0.47, "04.12.202516"),
    @(45995, 17, 0.251, "04.12.202517"),
    @(45995, 18, 0.091, "04.12.202518"),
    @(45995, 19, 0.043, "04.12.202519"),
    @(45995, 20, 0.048, "04.12.202520"),
    @(45995, 21, 0.048, "04.12.202521"),
    @(45995, 22, 0.016, "04.12.202522"),
    @(45995, 23, 0.011, "04.12.202523"),
    @(45995, 24, 0.109, "04.12.202524"),
    @(45996, 1, 0, "05.12.20251"),
    @(45996, 2, 0, "05.12.20252"),
    @(45996, 3, 0, "05.12.20253"),
    @(45996, 4, 0, "05.12.20254"),
    @(45996, 5, 0, "05.12.20255"),
    @(45996, 6, 0, "05.12.20256"),
    @(45996, 7, 0, "05.12.20257"),
    @(45996, 8, 0.035, "05.12.20258"),
    @(45996, 9, 0.08599999999999999, "05.12.20259"),
    @(45996, 10, 0.265, "05.12.202510"),
    @(45996, 11, 0.305, "05.12.202511"),
    @(45996, 12, 0.404, "05.12.202512"),
    @(45996, 13, 0.528, "05.12.202513"),
    @(45996, 14, 0.651, "05.12.202514"),
    @(45996, 15, 0.636, "05.12.202515"),
    @(45996, 16, 0.47, "05.12.202516"),
    @(45996, 17, 0.251, "05.12.202517"),
    @(45996, 18, 0.091, "05.12.202518"),
    @(45996, 19, 0.043, "05.12.202519"),
    @(45996, 20, 0.048, "05.12.202520"),
    @(45996, 21, 0.048, "05.12.202521"),
    @(45996, 22, 0.016, "05.12.202522"),
    @(45996, 23, 0.016, "05.12.202523"),
    @(45996, 24, 0.011, "05.12.202524"),
    @(45997, 1, 0, "06.12.20251"),
    @(45997, 2, 0, "06.12.20252"),
    @(45997, 3, 0, "06.12.20253"),
    @(45997, 4, 0, "06.12.20254"),
    @(45997, 5, 0, "06.12.20255"),
    @(45997, 6, 0, "06.12.20256"),
    @(45997, 7, 0, "06.12.20257"),
    @(45997, 8, 0.035, "06.12.20258"),
    @(45997, 9, 0.08599999999999999, "06.12.20259"),
    @(45997, 10, 0.265, "06.12.202510"),
    @(45997, 11, 0.488, "06.12.202511"),
    @(45997, 12, 0.706, "06.12.202512"),
    @(45997, 13, 1.021, "06.12.202513"),
    @(45997, 14, 0.721, "06.12.202514"),
    @(45997, 15, 0.636, "06.12.202515"),
    @(45997, 16, 0.402, "06.12.202516"),
    @(45997, 17, 0.251, "06.12.202517"),
    @(45997, 18, 0.091, "06.12.202518"),
    @(45997, 19, 0.043, "06.12.202519"),
    @(45997, 20, 0.048, "06.12.202520"),
    @(45997, 21, 0.048, "06.12.202521"),
    @(45997, 22, 0.016, "06.12.202522"),
    @(45997, 23, 0.016, "06.12.202523"),
    @(45997, 24, 0.016, "06.12.202524"),
    @(45998, 1, 0.013, "07.12.20251"),
    @(45998, 2, 0.013, "07.12.20252"),
    @(45998, 3, 0.013, "07.12.20253"),
    @(45998, 4, 0.013, "07.12.20254"),
    @(45998, 5, 0.013, "07.12.20255"),
    @(45998, 6, 0.013, "07.12.20256"),
    @(45998, 7, 0.013, "07.12.20257"),
    @(45998, 8, 0.04, "07.12.20258"),
    @(45998, 9, 0.091, "07.12.20259"),
    @(45998, 10, 0.27, "07.12.202510"),
    @(45998, 11, 0.31, "07.12.202511"),
    @(45998, 12, 0.409, "07.12.202512"),
    @(45998, 13, 0.533, "07.12.202513"),
    @(45998, 14, 0.59, "07.12.202514"),
    @(45998, 15, 0.469, "07.12.202515"),
    @(45998, 16, 0.336, "07.12.202516"),
    @(45998, 17, 0.251, "07.12.202517"),
    @(45998, 18, 0.091, "07.12.202518"),
    @(45998, 19, 0.043, "07.12.202519"),
    @(45998, 20, 0.048, "07.12.202520"),
    @(45998, 21, 0.048, "07.12.202521"),
    @(45998, 22, 0.016, "07.12.202522"),
    @(45998, 23, 0.016, "07.12.202523"),
    @(45998, 24, 0.016, "07.12.202524"),
    @(45999, 1, 0.013, "08.12.20251"),
    @(45999, 2, 0.013, "08.12.20252"),
    @(45999, 3, 0, "08.12.20253"),
    @(45999, 4, 0, "08.12.20254"),
    @(45999, 5, 0, "08.12.20255"),
    @(45999, 6, 0, "08.12.20256"),
    @(45999, 7, 0, "08.12.20257"),
    @(45999, 8, 0.035, "08.12.20258"),
    @(45999, 9, 0.08599999999999999, "08.12.20259"),
    @(45999, 10, 0.265, "08.12.202510"),
    @(45999, 11, 0.497, "08.12.202511"),
    @(45999, 12, 1.014, "08.12.202512"),
    @(45999, 13, 1.022, "08.12.202513"),
    @(45999, 14, 1.031, "08.12.202514"),
    @(45999, 15, 1.014, "08.12.202515"),
    @(45999, 16, 0.587, "08.12.202516"),
    @(45999, 17, 0.297, "08.12.202517"),
    @(45999, 18, 0.095, "08.12.202518"),
    @(45999, 19, 0.043, "08.12.202519"),
    @(45999, 20, 0.048, "08.12.202520"),
    @(45999, 21, 0.048, "08.12.202521"),
    @(45999, 22, 0.016, "08.12.202522"),
    @(45999, 23, 0.016, "08.12.202523"),
    @(45999, 24, 0.011, "08.12.202524"),
    @(46000, 1, 0, "09.12.20251"),
    @(46000, 2, 0, "09.12.20252"),
    @(46000, 3, 0, "09.12.20253"),
    @(46000, 4, 0, "09.12.20254"),
    @(46000, 5, 0, "09.12.20255"),
    @(46000, 6, 0, "09.12.20256"),
    @(46000, 7, 0.013, "09.12.20257"),
    @(46000, 8, 0.043, "09.12.20258"),
    @(46000, 9, 0.093, "09.12.20259"),
    @(46000, 10, 0.273, "09.12.202510"),
    @(46000, 11, 0.418, "09.12.202511"),
    @(46000, 12, 0.654, "09.12.202512"),
    @(46000, 13, 1.023, "09.12.202513")
)

$n = $data.Count
$arr = New-Object "object[,]" $n,4
for ($i = 0; $i -lt $n; $i++) {
    $arr[$i,0] = $data[$i][0]
    $arr[$i,1] = $data[$i][1]
    $arr[$i,2] = $data[$i][2]
    $arr[$i,3] = $data[$i][3]
}

$startRow = 2
$endRow = $startRow + $n - 1
$ws.Range("A" + $startRow + ":D" + $endRow).Value = $arr
